$wb = $excel.ActiveWorkbook

# --- Worksheet references ---
$wsBug  = $wb.Worksheets.Item("bug")
$wsNew  = $wb.Worksheets.Item("newfeature")

# --- Add the three new "newfeature" rows (42, 43, 44 in the "No." column -> sheet rows 43-45) ---
$wsNew.Range("C43").Value = "meetsdk"
$wsNew.Range("D43").Value = "android"
$wsNew.Range("E43").Value = "FFPlayer ogles2 render"
$wsNew.Range("F43").Value = "TBD"

$wsNew.Range("C44").Value = "meetsdk"
$wsNew.Range("D44").Value = "android"
$wsNew.Range("E44").Value = "XOPlayer support no-audio media"
$wsNew.Range("F44").Value = "TBD"

$wsNew.Range("C45").Value = "meetsdk"
$wsNew.Range("D45").Value = "android"
$wsNew.Range("E45").Value = "FFPlayer support render filter(ogles2)"
$wsNew.Range("F45").Value = "TBD"

# --- Switch the active sheet from "bug" to "newfeature" ---
$wsNew.Activate()

# --- Restore the view state: scroll position + selection on each sheet ---
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 4
$wsBug.Range("G53").Select()

$wsNew.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$wsNew.Range("F45").Select()
